$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.383.97"
$ws.Range("E2").Value = "  +5.75%  "
$ws.Range("D3").Value = "3.348.29"
$ws.Range("E3").Value = "  +2.34%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "413.94"
$ws.Range("E5").Value = "  +4.08%  "
$ws.Range("D6").Value = "112.43"
$ws.Range("E6").Value = "  +1.92%  "
$ws.Range("E7").Value = "  +4.30%  "
$ws.Range("D9").Value = "0.634"
$ws.Range("E9").Value = "  +1.17%  "
$ws.Range("D10").Value = "40.05"
$ws.Range("E10").Value = "  +1.89%  "
$ws.Range("D11").Value = "0.0988"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").Value = "3.880.71"
$ws.Range("E13").Value = "  +2.50%  "
$ws.Range("D14").Value = "8.53"
$ws.Range("E14").Value = "  +4.29%  "
$ws.Range("D15").Value = "19.36"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "3.352.40"
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").Value = "60.196.58"
$ws.Range("E18").Value = "  +5.81%  "
$ws.Range("D19").Value = "10.79"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("E21").Value = "  +4.22%  "
$ws.Range("D22").Value = "13.23"
$ws.Range("E22").Value = "  +2.04%  "
$ws.Range("D23").Value = "302.86"
$ws.Range("E23").Value = "  -1.11%  "
$ws.Range("D24").Value = "75.71"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").Value = "3.20"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").Value = "28.70"
$ws.Range("E26").Value = "  +1.42%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "4.49"
$ws.Range("E27").Value = "  +2.22%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.182"
$ws.Range("E28").Value = "  +6.59%  "
$ws.Range("D29").Value = "8.00"
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "7.48"
$ws.Range("E30").Value = "  +2.98%  "
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "2.69"
$ws.Range("E31").Value = "  +25.74%  "
$ws.Range("D32").Value = "0.116"
$ws.Range("E32").Value = "  +5.05%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "11.60"
$ws.Range("E33").Value = "  +4.91%  "
$ws.Range("B34").Value = "Dai"
$ws.Range("C34").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "39.85"
$ws.Range("E35").Value = "  +6.01%  "
$ws.Range("E36").Value = "  +5.67%  "
$ws.Range("D37").Value = "52.33"
$ws.Range("E37").Value = "  +1.43%  "
$ws.Range("D38").Value = "3.15"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("E40").Value = "  -3.93%  "
$ws.Range("D41").Value = "137.87"
$ws.Range("E41").Value = "  +1.95%  "
$ws.Range("E42").Value = "  +2.15%  "
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "0.286"
$ws.Range("E44").Value = "  +1.88%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "3.94"
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").Value = "16.95"
$ws.Range("E46").Value = "  -2.56%  "
$ws.Range("E47").Value = "  +8.89%  "
$ws.Range("D48").Value = "22.43"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("D49").Value = "2.208.46"
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("E50").Value = "  +1.83%  "
$ws.Range("D51").Value = "2.00"
$ws.Range("E51").Value = "  -1.63%  "